$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: classical-best-embeddings vs. classical-best-tfidf
$ws.Range("A2").Value = "classical-best-embed vs. classical-best-tfidf"
$ws.Range("C2").Value = 0.066
$ws.Range("D2").Value = 0.048
$ws.Range("E2").Value = 0.045
$ws.Range("I2").Value = 0.046
$ws.Range("J2").Value = 0.055

# Row 3: BERT-base vs. classical-best-tfidf (label unchanged)
$ws.Range("C3").Value = 0.07000000000000001
$ws.Range("D3").Value = 0.124
$ws.Range("E3").Value = 0.122
$ws.Range("F3").Value = 0.097
$ws.Range("G3").Value = 0.152
$ws.Range("H3").Value = 0.164
$ws.Range("I3").Value = 0.103
$ws.Range("J3").Value = 0.122

# Row 4: BERT-base vs. classical-best-embeddings
$ws.Range("A4").Value = "BERT-base vs. classical-best-embed"
$ws.Range("C4").Value = 0.004
$ws.Range("D4").Value = 0.076
$ws.Range("E4").Value = 0.077
$ws.Range("F4").Value = 0.07000000000000001
$ws.Range("G4").Value = 0.083
$ws.Range("H4").Value = 0.09
$ws.Range("I4").Value = 0.057
$ws.Range("J4").Value = 0.067

# Row 5: BERT-base-nli vs. classical-best-tfidf (label unchanged)
$ws.Range("B5").Value = 0.384
$ws.Range("C5").Value = 0.183
$ws.Range("E5").Value = 0.145
$ws.Range("G5").Value = 0.128
$ws.Range("H5").Value = 0.14
$ws.Range("J5").Value = 0.142

# Row 6: BERT-base-nli vs. classical-best-embeddings
$ws.Range("A6").Value = "BERT-base-nli vs. classical-best-embed"
$ws.Range("B6").Value = 0.384
$ws.Range("C6").Value = 0.117
$ws.Range("D6").Value = 0.103
$ws.Range("E6").Value = 0.1
$ws.Range("G6").Value = 0.059
$ws.Range("H6").Value = 0.066
$ws.Range("I6").Value = 0.1
$ws.Range("J6").Value = 0.08699999999999999

# Row 7: BERT-base-nli vs. BERT-base (label unchanged)
$ws.Range("B7").Value = 0.384
$ws.Range("C7").Value = 0.113
$ws.Range("D7").Value = 0.027
$ws.Range("E7").Value = 0.023
$ws.Range("F7").Value = 0.01
$ws.Range("G7").Value = -0.024
$ws.Range("H7").Value = -0.024
$ws.Range("I7").Value = 0.043
$ws.Range("J7").Value = 0.021
